# Generate Report for Handback
# Updates the "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" timestamp cells to reflect a freshly
# regenerated handback report.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
# "Latest HO Xliff Generate Date" for the first data row
$overview.Range("G2").Value = "2016-09-03 15:10:50"

# --- zh-cn sheet -------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
# "Correspond Handoff Datetime"
$zhcn.Range("H2").Value = "2016-09-03 15:10:46"
# "Correspond Handback DateTime"
$zhcn.Range("K2").Value = "2016-09-03 15:11:11"

# --- de-de sheet -------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
# "Correspond Handoff Datetime" (shares the same timestamp as Overview!G2)
$dede.Range("H2").Value = "2016-09-03 15:10:50"
# "Correspond Handback DateTime"
$dede.Range("K2").Value = "2016-09-03 15:11:19"
